# Apply the commit's changes:
#  - Rename "Paineis DARQ" -> "PAINEIS DARQ"
#  - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#  - Delete the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

$excel.DisplayAlerts = $false

$wsDarq = $wb.Worksheets.Item("Paineis DARQ")
$wsDarq.Name = "PAINEIS DARQ"

$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

$wsDesarq = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsDesarq.Delete()

$excel.DisplayAlerts = $true

# Keep "PAINEIS DARQ" as the selected/active tab, matching the source file.
$wsDarq.Activate()
